# This edit adds a new weekly price record at the top of the data table
# (row 229), pushing all the existing records (previously rows 229-314)
# down by one row (to 230-315). The sheet's used range therefore grows
# from A1:R314 to A1:R315.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 229, shifting rows
# 229-314 down to 230-315 (and extending the sheet dimension to R315).
$ws.Rows.Item(229).Insert()

# Populate the newly inserted row 229 with the new record.
$ws.Range("A229").Value = 11
$ws.Range("B229").Value = "Vega Monumental Concepción"
$ws.Range("C229").Value = "Bíobío"
$ws.Range("D229").Value = 45146
$ws.Range("E229").Value = 8
$ws.Range("F229").Value = 100112003
$ws.Range("G229").Value = "Ajo"
$ws.Range("H229").Value = "Chino"
$ws.Range("I229").Value = "Primera"
$ws.Range("J229").Value = 200
$ws.Range("K229").Value = 19000
$ws.Range("L229").Value = 19500
$ws.Range("M229").Value = 19250
$ws.Range("N229").Value = "`$/caja 10 kilos"
$ws.Range("O229").Value = "China"
$ws.Range("P229").Value = 1925
$ws.Range("Q229").Value = 10
$ws.Range("R229").Value = "Hortaliza"

# Match the date cell formatting used by the rest of column D.
$ws.Range("D229").NumberFormat = $ws.Range("D230").NumberFormat
